$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Starting layout:
#   P1: "23/04 "
#   P2 (underlined run props): TAB "-" "Equals"(spellcheck) " classe texto" + bookmark _GoBack
#
# Target layout:
#   P1: "23/04 " + bookmark _GoBack
#   P2: TAB "-" "Equals"(spellcheck) " classe texto"                     (new plain paragraph)
#   P3: (empty paragraph)                                                (new)
#   P4: "25/04"                                                          (new)
#   P5 (keeps the original underlined pPr): TAB "-" "Implementando"(gramcheck)
#       " a classe sala e reformulando a classe texto"
# ------------------------------------------------------------------

# Step 1: remove the existing _GoBack bookmark (it currently sits at the
# end of the original second paragraph; it will be re-created at the end
# of paragraph 1 in step 3).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Step 2: insert the three new paragraphs right after paragraph 1 and
# before the old "Equals classe texto" paragraph. Using raw OOXML via
# InsertXML lets us place the <w:proofErr> spell-check markers exactly
# like the original document did for "Equals".
$p1 = $d.Paragraphs(1).Range
$insertionPoint = $d.Range($p1.End - 1, $p1.End - 1)
$xmlNewParagraphs = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' `
    + '<w:p><w:r><w:tab/><w:t>-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Equals</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> classe texto</w:t></w:r></w:p>' `
    + '<w:p/>' `
    + '<w:p><w:r><w:t>25/04</w:t></w:r></w:p>' `
    + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$insertionPoint.InsertXML($xmlNewParagraphs)

# Step 3: re-create the _GoBack bookmark at the end of paragraph 1
# (right after "23/04 ", before its paragraph mark). Adding a bookmark
# directly at that exact boundary collapses to the wrong range in this
# runtime, so we temporarily pad the paragraph with one extra character,
# anchor the bookmark there, and then remove the padding again.
$p1 = $d.Paragraphs(1).Range
$tailPos = $p1.End - 1
$pad = $d.Range($tailPos, $tailPos)
$pad.InsertAfter("X")
$bookmarkTarget = $d.Range($tailPos, $tailPos)
[void]$d.Bookmarks.Add("_GoBack", $bookmarkTarget)
$padRange = $d.Range($tailPos, $tailPos + 1)
[void]$padRange.Delete()

# Step 4: replace the content of the final paragraph (the original
# "Equals classe texto" paragraph, which still carries the underlined
# paragraph formatting) with the new "Implementando..." text, swapping
# the spell-check proofErr markers for grammar-check ones.
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaIndex).Range
$lastParaContent = $d.Range($lastPara.Start, $lastPara.End - 1)
$xmlReplacement = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' `
    + '<w:r><w:tab/><w:t>-</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Implementando</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> a classe sala e reformulando a classe texto</w:t></w:r>' `
    + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$lastParaContent.InsertXML($xmlReplacement)
